$wb = $excel.ActiveWorkbook

# Update the "想去人数" (column F) values on the "展览" sheet
$wsExhibition = $wb.Worksheets.Item(1)
$wsExhibition.Range("F3").Value = 0
$wsExhibition.Range("F4").Value = 0
$wsExhibition.Range("F5").Value = 450
$wsExhibition.Range("F6").Value = 155
$wsExhibition.Range("F9").Value = 0
$wsExhibition.Range("F10").Value = 0
$wsExhibition.Range("F17").Value = 0
$wsExhibition.Range("F18").Value = 0
$wsExhibition.Range("F19").Value = 5056
$wsExhibition.Range("F22").Value = 488
$wsExhibition.Range("F23").Value = 0

# Update the "想去人数" (column F) values on the "全部类型" sheet
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value = 6905
$wsAll.Range("F4").Value = 0
$wsAll.Range("F6").Value = 0
$wsAll.Range("F8").Value = 0
$wsAll.Range("F12").Value = 0
$wsAll.Range("F13").Value = 0
$wsAll.Range("F14").Value = 0
$wsAll.Range("F16").Value = 397
$wsAll.Range("F17").Value = 0
$wsAll.Range("F18").Value = 0
$wsAll.Range("F20").Value = 0
$wsAll.Range("F21").Value = 46
$wsAll.Range("F22").Value = 105
$wsAll.Range("F24").Value = 0
$wsAll.Range("F26").Value = 204
